$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Evcard: add a new trip logged at the end of 2018 (row 89)
# ---------------------------------------------------------------------
$evcard = $wb.Worksheets.Item("Evcard")

# Copy the formatting from the row above (row 88) down onto the new
# row so the new cells pick up the same styles (date format / borders)
# as the rest of the table, instead of Excel inventing a brand-new
# number format.
$evcard.Range("A88:B88").Copy()
$evcard.Range("A89:B89").PasteSpecial(-4122)

$evcard.Cells.Item(89, 1).Value = 43465
$evcard.Cells.Item(89, 2).Value = 3
$evcard.Cells.Item(89, 3).Formula = "=SUM(B82:B89)"

[void]$evcard.Activate()
[void]$evcard.Range("F80").Select()

# ---------------------------------------------------------------------
# Ponycar: total the last batch of rows (row 51) with a SUM formula
# ---------------------------------------------------------------------
$ponycar = $wb.Worksheets.Item("Ponycar")
$ponycar.Cells.Item(51, 3).Formula = "=SUM(B48:B51)"

[void]$ponycar.Activate()
[void]$ponycar.Range("C52").Select()

# ---------------------------------------------------------------------
# SUM: the Evcard total (whole-column SUM) needs to pick up the new
# B89 value added above; re-apply the formula so it recalculates.
# ---------------------------------------------------------------------
$sum = $wb.Worksheets.Item("SUM")
$sum.Cells.Item(2, 2).Formula = "=SUM(Evcard!B:B)"
